# Add the new project-plan rows ("mongo db and docs") below the existing
# table (row 11 left blank, matching the author's spreadsheet layout).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Creation of mongodb database"
$ws.Range("B12").Value = 2

$ws.Range("A13").Value = "creation of nodejs web api"
$ws.Range("B13").Value = 2

$ws.Range("A14").Value = "pulling of datafrom mongo db to google maps application"
$ws.Range("B14").Value = 1

$ws.Range("A15").Value = "creation of tags and putting them into mongodb"
$ws.Range("B15").Value = 1

$ws.Range("A16").Value = "ui/css of web app"
$ws.Range("B16").Value = 1

# Leave the selection where the author ended up after entering the last value.
$ws.Range("B17").Select()
